$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("测试总况")

# 测试日期 (test date): 2016-09-07 11:24 AM -> 2016-09-07 11:34 AM
$ws.Range("C6").Value = "2016-09-07 11:34 AM"

# 测试耗时 (test duration): 21 -> 30
# Force text so the digits-only string isn't auto-converted to a number.
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "30"

# 内存占用均值 (avg memory usage): 23% -> 22%
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22%"

# 内存占用峰值 (peak memory usage): 78070KB -> 77021KB
$ws.Range("D12").Value = "77021KB"

# CPU占用均值 (avg CPU usage): 36% -> 63%
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "63%"

# CPU占用峰值 (peak CPU usage): 53% -> 71%
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "71%"
